$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply base formatting (border + left + vcenter, matches existing style of row 1) ---
# to the full new extent A1:E10 so every cell (old + new) carries the same format.
$ws.Range("A1").Copy()
$ws.Range("A1:E10").PasteSpecial(-4122)

# --- Header row ---
$ws.Range("A1").Value = " 序号"
$ws.Range("B1").Value = "PIN"
$ws.Range("C1").Value = "方向"
$ws.Range("D1").Value = "说明"
$ws.Range("E1").Value = "外设映射"

# --- Data rows (PIN, direction, description, peripheral mapping) ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "P10"
$ws.Range("C2").Value = "输出"
$ws.Range("D2").Value = "PWM输出(控制亮度)"
$ws.Range("E2").Value = "PWM3"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "P03"
$ws.Range("C3").Value = "输出"
$ws.Range("D3").Value = "PWM输出(控制冷暖)"
$ws.Range("E3").Value = "PWM0"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "P06"
$ws.Range("C4").Value = "输出"
$ws.Range("D4").Value = "PWM输出(控制R)"
$ws.Range("E4").Value = "PWM01"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "P00"
$ws.Range("C5").Value = "输出"
$ws.Range("D5").Value = "PWM输出(控制G)"
$ws.Range("E5").Value = "PWM1"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "P01"
$ws.Range("C6").Value = "输出"
$ws.Range("D6").Value = "PWM输出(控制B)"
$ws.Range("E6").Value = "PWM11"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "P02"
$ws.Range("C7").Value = "模拟输入"
$ws.Range("D7").Value = "雷达中频输入"
$ws.Range("E7").Value = "AN2"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "P22"
$ws.Range("C8").Value = "模拟输入"
$ws.Range("D8").Value = "光敏电阻输入"
$ws.Range("E8").Value = "AN10"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "P20"
$ws.Range("C9").Value = "TXD"
$ws.Range("D9").Value = "MCU串口发"
$ws.Range("E9").Value = "TXD"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "P04"
$ws.Range("C10").Value = "RXD"
$ws.Range("D10").Value = "MCU串口收"
$ws.Range("E10").Value = "RXD"

# --- E9:E10 use a variant style: same border + vcenter, but without the
#     horizontal="left" alignment override (matches the new 4th cellXfs entry). ---
$ws.Range("E9:E10").HorizontalAlignment = 1

# --- Clear the stale single-cell selection left over from the old C7 selection ---
$ws.Range("A1").Select()
